# Update the "as_of_utc" timestamp column (AA) for every data row
# (rows 2-26) on both the "Главные" and "Линейные" sheets.
$wb = $excel.ActiveWorkbook

$newTimestamp = "2025-12-19 07:03:24"
$sheetNames = @("Главные", "Линейные")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($row = 2; $row -le 26; $row++) {
        $ws.Cells.Item($row, 27).Value = $newTimestamp
    }
}
